$wb = $excel.ActiveWorkbook

# --- AccountDetails sheet: update address2 (H2) and mobileNumber (M2) ---
$ws = $wb.Worksheets.Item("AccountDetails")
$ws.Range("H2").Value = "Karachi 2020"
$ws.Range("M2").Value = "536691048"

# Widen columns I:K (country/state/city) to match the reviewed layout.
$ws.Columns.Item(9).ColumnWidth = 16.26
$ws.Columns.Item(10).ColumnWidth = 13.76
$ws.Columns.Item(11).ColumnWidth = 10.43

# Make AccountDetails the active sheet/tab and move the selection to M3.
$ws.Activate()
$ws.Range("M3").Select()

# --- PaymentDetails sheet: no longer the active tab, selection stays E5 ---
$wsPay = $wb.Worksheets.Item("PaymentDetails")
$wsPay.Range("E5").Select()

# Re-activate AccountDetails last so it is the sheet shown/selected on open.
$ws.Activate()
